# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.621.16"
$ws.Range("E2").Value = "'  +1.50%  "
$ws.Range("D3").Value = "'2.473.04"
$ws.Range("E3").Value = "'  +1.47%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("D5").Value = "'575.91"
$ws.Range("E5").Value = "'  +1.54%  "
$ws.Range("D6").Value = "'148.78"
$ws.Range("E6").Value = "'  +2.35%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "'  +1.59%  "
$ws.Range("D9").Value = "'2.466.69"
$ws.Range("E9").Value = "'  +1.15%  "
$ws.Range("E10").Value = "'  +0.70%  "
$ws.Range("D12").Value = "'5.30"
$ws.Range("E12").Value = "'  +0.81%  "
$ws.Range("D13").Value = "'0.358"
$ws.Range("E13").Value = "'  +1.22%  "
$ws.Range("D14").Value = "'27.24"
$ws.Range("E14").Value = "'  +1.51%  "
$ws.Range("D15").Value = "'0.0000181"
$ws.Range("E15").Value = "'  -1.18%  "
$ws.Range("D16").Value = "'2.919.71"
$ws.Range("E16").Value = "'  +1.46%  "
$ws.Range("D17").Value = "'63.565.78"
$ws.Range("E17").Value = "'  +1.53%  "
$ws.Range("D18").Value = "'2.479.43"
$ws.Range("E18").Value = "'  +1.87%  "
$ws.Range("D19").Value = "'11.47"
$ws.Range("E19").Value = "'  +2.15%  "
$ws.Range("E20").Value = "'  +7.18%  "
$ws.Range("D21").Value = "'330.55"
$ws.Range("E21").Value = "'  +2.22%  "
$ws.Range("D22").Value = "'4.22"
$ws.Range("E22").Value = "'  +1.30%  "
$ws.Range("E23").Value = "'  +18.43%  "
$ws.Range("E24").Value = "'  +0.19%  "
$ws.Range("D25").Value = "'65.97"
$ws.Range("E25").Value = "'  -1.96%  "
$ws.Range("D26").Value = "'629.89"
$ws.Range("E26").Value = "'  +11.28%  "
$ws.Range("D27").Value = "'9.21"
$ws.Range("E27").Value = "'  +5.81%  "
$ws.Range("D28").Value = "'0.0000105"
$ws.Range("E28").Value = "'  +3.51%  "
$ws.Range("E29").Value = "'  +5.96%  "
$ws.Range("D30").Value = "'2.603.37"
$ws.Range("B31").Value = "'Binance-PegBSC-USD"
$ws.Range("C31").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "'  +0.07%  "
$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.41"
$ws.Range("E32").Value = "'  +0.15%  "
$ws.Range("E33").Value = "'  -2.35%  "
$ws.Range("E34").Value = "'  +1.59%  "
$ws.Range("D35").Value = "'5.26"
$ws.Range("E35").Value = "'  +8.21%  "
$ws.Range("E36").Value = "'  +0.42%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "'  +0.00%  "
$ws.Range("D38").Value = "'0.383"
$ws.Range("E38").Value = "'  +0.11%  "
$ws.Range("D39").Value = "'5.50"
$ws.Range("E39").Value = "'  +0.96%  "
$ws.Range("D40").Value = "'18.88"
$ws.Range("E40").Value = "'  +0.43%  "
$ws.Range("D41").Value = "'2.75"
$ws.Range("E41").Value = "'  +13.70%  "
$ws.Range("B42").Value = "'Monero"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'147.25"
$ws.Range("E42").Value = "'  -0.69%  "
$ws.Range("B43").Value = "'Stacks"
$ws.Range("C43").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.82"
$ws.Range("E43").Value = "'  -0.43%  "
$ws.Range("E44").Value = "'  -0.45%  "
$ws.Range("D45").Value = "'150.98"
$ws.Range("E45").Value = "'  +1.72%  "
$ws.Range("D46").Value = "'3.78"
$ws.Range("E46").Value = "'  +2.98%  "
$ws.Range("D47").Value = "'21.56"
$ws.Range("E47").Value = "'  +5.21%  "
$ws.Range("D48").Value = "'0.0542"
$ws.Range("E48").Value = "'  +0.71%  "
$ws.Range("D49").Value = "'0.606"
$ws.Range("E49").Value = "'  +0.98%  "
$ws.Range("D50").Value = "'0.0236"
$ws.Range("E50").Value = "'  +2.16%  "
$ws.Range("E51").Value = "'  -0.79%  "
